# Insert a new weekly price-report row for "Vega Modelo de Temuco - Achicoria".
# The sheet's data rows shift down by one starting at row 24 (a new sample
# dated 2022-07-22 / serial 44771 is inserted ahead of the existing series),
# so every subsequent row's content effectively moves down by one position
# and a brand-new row is written at row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 24..55 down to 25..56, carrying formatting (date style, etc.)
# along with them, and leave a blank row 24 ready to be filled in.
$ws.Rows(24).Insert()

# Populate the newly inserted row 24 with the new observation.
$ws.Cells.Item(24, 1).Value = 10
$ws.Cells.Item(24, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(24, 3).Value = "La Araucanía"
$ws.Cells.Item(24, 4).Value = 44771
$ws.Cells.Item(24, 5).Value = 9
$ws.Cells.Item(24, 6).Value = 100112010
$ws.Cells.Item(24, 7).Value = "Achicoria"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 9000
$ws.Cells.Item(24, 12).Value = 9000
$ws.Cells.Item(24, 13).Value = 9000
$ws.Cells.Item(24, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 500
$ws.Cells.Item(24, 17).Value = 18
$ws.Cells.Item(24, 18).Value = "Hortaliza"
